$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.965.73'
$ws.Range("E2").Value = '  -1.08%  '

$ws.Range("D3").Value = '3.374.25'
$ws.Range("E3").Value = '  -0.25%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.44'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.23%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").Value = '3.373.20'
$ws.Range("E8").Value = '  -0.26%  '

$ws.Range("E9").Value = '  -1.35%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.60'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.74%  '

$ws.Range("E11").Value = '  -3.29%  '

$ws.Range("E12").Value = '  -2.61%  '

$ws.Range("D13").Value = '3.949.09'
$ws.Range("E13").Value = '  -0.27%  '

$ws.Range("E14").Value = '  +0.53%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.63'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.01%  '

$ws.Range("D16").Value = '3.369.74'
$ws.Range("E16").Value = '  -0.52%  '

$ws.Range("E17").Value = '  -3.24%  '

$ws.Range("D18").Value = '61.098.20'
$ws.Range("E18").Value = '  -0.99%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.00%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.72'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.32'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '373.79'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.38%  '

$ws.Range("D23").Value = '3.517.01'
$ws.Range("E23").Value = '  -0.22%  '

$ws.Range("E24").Value = '  -2.59%  '

$ws.Range("E25").Value = '  -0.10%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '71.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.31%  '

$ws.Range("E27").Value = '  -1.17%  '

$ws.Range("E28").Value = '  -6.15%  '

$ws.Range("E29").Value = '  +10.93%  '

$ws.Range("E30").Value = '  -0.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.36'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.13'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.83%  '

$ws.Range("E34").Value = '  -0.07%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.28'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.12'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.35%  '

$ws.Range("E37").Value = '  -1.20%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.78'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.92%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '164.55'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0756'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.64%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.01%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.774'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.80%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '24.72'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.67'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.77%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.30'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.31%  '

$ws.Range("E46").Value = '  -4.74%  '

$ws.Range("D47").Value = '2.534.37'
$ws.Range("E47").Value = '  +8.52%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.76'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.49%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.81'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.06%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.41'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.03%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0257'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.47%  '
